$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Remember the date-column number format used on the original data row
# before we start overwriting cells.
$dateFmt = $ws.Cells.Item(5, 1).NumberFormat

# The existing row 5 (the original single trade) moves down to row 10.
# Write it to row 10 first (order doesn't matter since we don't use Insert).
$ws.Cells.Item(10, 1).Value = 46062
$ws.Cells.Item(10, 1).NumberFormat = $dateFmt
$ws.Cells.Item(10, 2).Value = "NSE"
$ws.Cells.Item(10, 3).Value = "Buy"
$ws.Cells.Item(10, 4).Value = 100
$ws.Cells.Item(10, 5).Value = 36.7
$ws.Cells.Item(10, 6).Value = 3695.58
$ws.Cells.Item(10, 7).Value = "CN#252611665409"
$ws.Cells.Item(10, 8).Value = 3.67
$ws.Cells.Item(10, 9).Value = 21.91
$ws.Cells.Item(10, 10).Formula = "=Index!`$C`$2"

# New trade rows 5-9 inserted above the old row (now at row 10).
# Row 5
$ws.Cells.Item(5, 1).Value = 46059
$ws.Cells.Item(5, 2).Value = "NSE"
$ws.Cells.Item(5, 3).Value = "Buy"
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 37.32
$ws.Cells.Item(5, 6).Value = 37.32
$ws.Cells.Item(5, 7).Value = "~"
$ws.Cells.Item(5, 8).Clear()
$ws.Cells.Item(5, 9).Clear()
$ws.Cells.Item(5, 10).Formula = "=Index!`$C`$2"

# Row 6
$ws.Cells.Item(6, 1).Value = 46059
$ws.Cells.Item(6, 2).Value = "NSE"
$ws.Cells.Item(6, 3).Value = "Buy"
$ws.Cells.Item(6, 4).Value = 99
$ws.Cells.Item(6, 5).Value = 37.31
$ws.Cells.Item(6, 6).Value = 3693.69
$ws.Cells.Item(6, 7).Value = "~"
$ws.Cells.Item(6, 10).Formula = "=Index!`$C`$2"

# Row 7
$ws.Cells.Item(7, 1).Value = 46050
$ws.Cells.Item(7, 2).Value = "NSE"
$ws.Cells.Item(7, 3).Value = "Buy"
$ws.Cells.Item(7, 4).Value = 100
$ws.Cells.Item(7, 5).Value = 38.21
$ws.Cells.Item(7, 6).Value = 3821
$ws.Cells.Item(7, 7).Value = "~"
$ws.Cells.Item(7, 10).Formula = "=Index!`$C`$2"

# Row 8
$ws.Cells.Item(8, 1).Value = 46049
$ws.Cells.Item(8, 2).Value = "NSE"
$ws.Cells.Item(8, 3).Value = "Buy"
$ws.Cells.Item(8, 4).Value = 100
$ws.Cells.Item(8, 5).Value = 35.88
$ws.Cells.Item(8, 6).Value = 3588
$ws.Cells.Item(8, 7).Value = "~"
$ws.Cells.Item(8, 10).Formula = "=Index!`$C`$2"

# Row 9
$ws.Cells.Item(9, 1).Value = 46044
$ws.Cells.Item(9, 2).Value = "NSE"
$ws.Cells.Item(9, 3).Value = "Buy"
$ws.Cells.Item(9, 4).Value = 250
$ws.Cells.Item(9, 5).Value = 39.38
$ws.Cells.Item(9, 6).Value = 9845
$ws.Cells.Item(9, 7).Value = "~"
$ws.Cells.Item(9, 10).Formula = "=Index!`$C`$2"

# Column A (date column) uses a custom date number format on every data row;
# match that for the new rows too.
$ws.Cells.Item(5, 1).NumberFormat = $dateFmt
$ws.Cells.Item(6, 1).NumberFormat = $dateFmt
$ws.Cells.Item(7, 1).NumberFormat = $dateFmt
$ws.Cells.Item(8, 1).NumberFormat = $dateFmt
$ws.Cells.Item(9, 1).NumberFormat = $dateFmt
